$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (modelo_2012) and C (modelo_2016), rows 2-32
$bValues = @(0,0,0,15,20,60,95,100,50,5,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
$cValues = @(0,0,0,0,0,15,20,35,75,95,100,90,60,5,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}

# Update the active selection from M6 to C1
$ws.Range("C1").Select()
